# "Add files via upload" - update sample station labels in the DY79 sheet.
# Column A (rows 2-49) previously held ad-hoc barcode-style sample names
# (BC01, BC02, ...). They are renamed sequentially to Station1..Station48,
# and the sheet view/selection + column A width are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DY79")

for ($i = 2; $i -le 49; $i++) {
    $stationNum = $i - 1
    $ws.Cells.Item($i, 1).Value = "Station$stationNum"
}

# Column A now holds longer labels ("Station47", etc.) than the old BC
# codes, so widen it to fit (matches the custom width seen after editing).
$ws.Columns.Item(1).AutoFit()

# Reflect the cursor position left after the edits were made.
$ws.Range("B8").Select() | Out-Null
